$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.988.10"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").Value = "1.880.39"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("D4").Value = "'1.012"
$ws.Range("E4").Value = "  +0.85%  "
$ws.Range("D5").Value = "'336.09"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").Value = "'1.009"
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("D7").Value = "'0.4774"
$ws.Range("E7").Value = "  +1.32%  "
$ws.Range("D8").Value = "'0.3947"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").Value = "'47.22"
$ws.Range("E9").Value = "  -1.62%  "
$ws.Range("D10").Value = "'0.08024"
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("D11").Value = "'1.019"
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").Value = "'21.90"
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("D13").Value = "1.885.15"
$ws.Range("E13").Value = "  +0.75%  "
$ws.Range("D14").Value = "'6.051"
$ws.Range("E14").Value = "  +1.70%  "
$ws.Range("D15").Value = "'7.205"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("E16").Value = "  +0.92%  "
$ws.Range("D17").Value = "'88.76"
$ws.Range("E17").Value = "  +2.29%  "
$ws.Range("D18").Value = "'0.06725"
$ws.Range("E18").Value = "  +1.03%  "
$ws.Range("D19").Value = "'0.00001051"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").Value = "'17.08"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").Value = "'1.009"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").Value = "27.972.30"
$ws.Range("E22").Value = "  +1.09%  "
$ws.Range("D23").Value = "'5.514"
$ws.Range("E23").Value = "  +0.28%  "
$ws.Range("D24").Value = "'11.01"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").Value = "'2.342"
$ws.Range("D26").Value = "2.105.22"
$ws.Range("E26").Value = "  +0.67%  "
$ws.Range("D27").Value = "'158.46"
$ws.Range("E27").Value = "  -0.25%  "
$ws.Range("D28").Value = "'19.91"
$ws.Range("E28").Value = "  -1.58%  "
$ws.Range("D29").Value = "'2.108"
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("D30").Value = "'5.498"
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("D31").Value = "'121.68"
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("D32").Value = "'0.9808"
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("D33").Value = "'0.09569"
$ws.Range("E33").Value = "  +0.61%  "
$ws.Range("D34").Value = "'3.632"
$ws.Range("E34").Value = "  +1.15%  "
$ws.Range("D35").Value = "'5.339"
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("D36").Value = "'1.357"
$ws.Range("E36").Value = "  -6.28%  "
$ws.Range("D37").Value = "'0.06081"
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("D38").Value = "'0.02247"
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("D39").Value = "'1.207"
$ws.Range("E39").Value = "  -1.97%  "
$ws.Range("D40").Value = "'8.199"
$ws.Range("E40").Value = "  +1.17%  "
$ws.Range("D41").Value = "'1.009"
$ws.Range("E41").Value = "  +0.66%  "
$ws.Range("D42").Value = "'0.5996"
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("D43").Value = "'0.1897"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "'10.36"
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("D45").Value = "'1.261"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("D46").Value = "'0.5686"
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("D47").Value = "'12.28"
$ws.Range("E47").Value = "  +0.38%  "
$ws.Range("D48").Value = "'1.934"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("E49").Value = "  -1.11%  "
$ws.Range("D50").Value = "'0.06810"
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("D51").Value = "'112.73"
$ws.Range("E51").Value = "  -1.38%  "

Write-Host "Applied cryptos update"